# This document's edit is a cyclic rotation of paragraph text contents:
# the paragraph *structure* (styles, runs, formatting, bookmarks, etc.)
# stays exactly where it is; only the text carried by specific <w:t>
# runs moves to a different (existing) paragraph slot.
#
# Two independent cycles of text are rotated:
#
#   Cycle 1 (Portuguese paragraphs):
#     A (Complementar a formacao ...)      -> slot that held B
#     B (A definir de acordo ...)          -> slot that held C
#     C (O conteudo desta disciplina ...)  -> slot that held D
#     D (Esta disciplina devera ...)       -> slot that held E
#     E (Media ponderada ...)              -> slot that held F
#     F (A recuperacao sera composta ...)  -> slot that held G
#     G (Livros, artigos ...)              -> slot that held H
#     H (5817650 - Erica Leonor Romao)     -> slot that held A
#
#   Cycle 2 (English/italic paragraphs):
#     X (Complement the training ...)      -> slot that held Y
#     Y (To be defined according ...)      -> slot that held X
#
# We implement this with Find/Replace, routing every source text through a
# unique temporary placeholder first so that no replacement can accidentally
# match text that was already moved into place earlier in the script.

$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replacement not found for: $old"
    }
}

# --- Original texts ---
$A = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$B = "A definir de acordo com o tópico programado"
$C = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$D = "Esta disciplina deverá conter no mínimo duas avaliações denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa. As avalições podem ser: seminários, trabalhos, projetos ou outra forma de avaliação definida pelo professor. Sendo necessário no mínimo uma avaliação na forma de prova escrita."
$E = "Média ponderada das avaliações (M)."
$F = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$G = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Meio Ambiente."
$H = "5817650 - Érica Leonor Romão"

$X = "Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art"
$Y = "To be defined according to the scheduled topic"

# --- Step 1: move every source text to a unique, unambiguous placeholder ---
Replace-Exact $A "@@PLACEHOLDER_A@@"
Replace-Exact $B "@@PLACEHOLDER_B@@"
Replace-Exact $C "@@PLACEHOLDER_C@@"
Replace-Exact $D "@@PLACEHOLDER_D@@"
Replace-Exact $E "@@PLACEHOLDER_E@@"
Replace-Exact $F "@@PLACEHOLDER_F@@"
Replace-Exact $G "@@PLACEHOLDER_G@@"
Replace-Exact $H "@@PLACEHOLDER_H@@"

Replace-Exact $X "@@PLACEHOLDER_X@@"
Replace-Exact $Y "@@PLACEHOLDER_Y@@"

# --- Step 2: place final texts from the placeholders, per the rotation ---
# The slot that used to hold A must now show B's text, B's slot shows C's
# text, ..., and H's slot (the last in the cycle) wraps around to show A's
# original text.
Replace-Exact "@@PLACEHOLDER_A@@" $B
Replace-Exact "@@PLACEHOLDER_B@@" $C
Replace-Exact "@@PLACEHOLDER_C@@" $D
Replace-Exact "@@PLACEHOLDER_D@@" $E
Replace-Exact "@@PLACEHOLDER_E@@" $F
Replace-Exact "@@PLACEHOLDER_F@@" $G
Replace-Exact "@@PLACEHOLDER_G@@" $H
Replace-Exact "@@PLACEHOLDER_H@@" $A

Replace-Exact "@@PLACEHOLDER_X@@" $Y
Replace-Exact "@@PLACEHOLDER_Y@@" $X

Write-Output "Done"
